# "POP read excel data" - update the sample/test fixture:
#  - rename the header cell from "Słowo" to "Fraza"
#  - drop the bottom border on the header cell (left/right/top stay thin/black)
#  - move the active selection to A2 (was C7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header text change (A1 keeps its existing bold/yellow-fill/bordered look,
#    just the caption changes)
$ws.Range("A1").Value = "Fraza"

# 2) Remove only the bottom edge of A1's border, keeping left/right/top thin+black
$ws.Range("A1").Borders.Item(9).LineStyle = -4142

# 3) Move the selection/active cell to A2
$ws.Range("A2").Select()
